$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.566.85"
$ws.Range("D3").Value = "2.609.48"
$ws.Range("E3").Value = "  +0.83%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'538.00"
$ws.Range("E5").Value = "  +2.92%  "
$ws.Range("D6").Value = "'141.66"
$ws.Range("E6").Value = "  +1.73%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +0.38%  "
$ws.Range("E9").Value = "  -0.26%  "
$ws.Range("E10").Value = "  +1.24%  "
$ws.Range("E11").Value = "  +1.40%  "
$ws.Range("E12").Value = "  -1.08%  "
$ws.Range("D13").Value = "3.067.13"
$ws.Range("E13").Value = "  +0.65%  "
$ws.Range("D14").Value = "59.497.01"
$ws.Range("E14").Value = "  +0.85%  "
$ws.Range("D15").Value = "'20.72"
$ws.Range("E15").Value = "  +1.21%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "'0.0000133"
$ws.Range("E16").Value = "  +0.55%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.590.88"
$ws.Range("E17").Value = "  -0.83%  "
$ws.Range("D18").Value = "'340.74"
$ws.Range("E19").Value = "  +1.39%  "
$ws.Range("D21").Value = "'6.32"
$ws.Range("E21").Value = "  -2.00%  "
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").Value = "'67.23"
$ws.Range("E23").Value = "  +1.23%  "
$ws.Range("E24").Value = "  +1.05%  "
$ws.Range("E25").Value = "  -1.41%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").Value = "'7.22"
$ws.Range("E27").Value = "  +2.54%  "
$ws.Range("D28").Value = "0.0₃0744"
$ws.Range("E28").Value = "  +3.00%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("E30").Value = "  +5.31%  "
$ws.Range("D31").Value = "'5.82"
$ws.Range("E31").Value = "  -1.22%  "
$ws.Range("D32").Value = "'18.76"
$ws.Range("E32").Value = "  +0.25%  "
$ws.Range("D33").Value = "'150.62"
$ws.Range("E33").Value = "  +0.86%  "
$ws.Range("E34").Value = "  +0.57%  "
$ws.Range("E35").Value = "  +0.68%  "
$ws.Range("D36").Value = "'0.846"
$ws.Range("E36").Value = "  +3.84%  "
$ws.Range("E37").Value = "  -0.86%  "
$ws.Range("D38").Value = "'0.824"
$ws.Range("E38").Value = "  +0.07%  "
$ws.Range("E39").Value = "  +0.39%  "
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("D41").Value = "'275.40"
$ws.Range("E41").Value = "  +1.30%  "
$ws.Range("E42").Value = "  +0.16%  "
$ws.Range("D43").Value = "'10.74"
$ws.Range("E43").Value = "  -0.48%  "
$ws.Range("E44").Value = "  -0.33%  "
$ws.Range("D45").Value = "'0.0523"
$ws.Range("E45").Value = "  +1.60%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "1.938.66"
$ws.Range("E46").Value = "  -1.49%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "'0.0223"
$ws.Range("E47").Value = "  +0.55%  "
$ws.Range("D48").Value = "'18.44"
$ws.Range("E48").Value = "  +2.54%  "
$ws.Range("E49").Value = "  +0.80%  "
$ws.Range("D50").Value = "'111.10"
$ws.Range("E50").Value = "  -2.31%  "
$ws.Range("E51").Value = "  +1.82%  "
